$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Refresh the "time_taken" query timestamps in column F on the "data" sheet
# (re-run happened ~40 minutes later than the original export).
$ws.Range("F2").Value = "2021-10-05 14:20:29.029406"
$ws.Range("F3").Value = "2021-10-05 14:20:29.029414"
$ws.Range("F4").Value = "2021-10-05 14:20:29.029417"
$ws.Range("F5").Value = "2021-10-05 14:20:29.029420"
$ws.Range("F6").Value = "2021-10-05 14:20:29.029423"
$ws.Range("F7").Value = "2021-10-05 14:20:29.029425"
$ws.Range("F8").Value = "2021-10-05 14:20:29.029428"
$ws.Range("F9").Value = "2021-10-05 14:20:29.029431"
$ws.Range("F10").Value = "2021-10-05 14:20:29.029434"
$ws.Range("F11").Value = "2021-10-05 14:20:29.029436"
$ws.Range("F12").Value = "2021-10-05 14:20:29.029439"
$ws.Range("F13").Value = "2021-10-05 14:20:29.029441"
$ws.Range("F14").Value = "2021-10-05 14:20:29.029444"
$ws.Range("F15").Value = "2021-10-05 14:20:29.029446"
$ws.Range("F16").Value = "2021-10-05 14:20:29.029449"
$ws.Range("F17").Value = "2021-10-05 14:20:29.029451"
$ws.Range("F18").Value = "2021-10-05 14:20:29.029454"
$ws.Range("F19").Value = "2021-10-05 14:20:29.029457"
$ws.Range("F20").Value = "2021-10-05 14:20:29.029459"
$ws.Range("F21").Value = "2021-10-05 14:20:29.029462"
$ws.Range("F22").Value = "2021-10-05 14:20:29.029464"

# Add the new "metadata" tab right after "data"
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "metadata"

# Reuse the bold / bordered / centered header style from the "data" sheet's
# header row (B1:F1) for the metadata header row (B1:G1).
$ws.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Reuse the index-column style (A2) from the "data" sheet for the metadata
# row index cell.
$ws.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").Value = "Gene therapy clinical trials"
$newSheet.Range("C2").Value = 412

# "0.7" is a version string, not a number - force text storage (like the
# source data's inline strings) then drop the temporary number-format style
# so the cell ends up with the plain/default format.
$d2 = $newSheet.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "0.7"
$d2.ClearFormats()

$newSheet.Range("E2").Value = "2019-06-20T15:11:44.609881Z"
$newSheet.Range("F2").Value = "2021-10-05 14:20:29.025669"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/412/?format=json"
